$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh the rolled-up metrics after trade #47 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.59   # Current Capital
$summary.Range("B4").Value = 0.59      # Total P&L $
$summary.Range("B5").Value = 0.25      # Total P&L %
$summary.Range("B6").Value = 47        # Total Trades
$summary.Range("B8").Value = 24        # Losing Trades
$summary.Range("B9").Value = 29.79     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) picks up the same refresh.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.59     # Capital
$status.Range("D4").Value = 47         # Trades
$status.Range("E4").Value = 0.59       # P&L $
$status.Range("F4").Value = 0.59       # P&L %
$status.Range("G4").Value = 29.79      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly-closed trade (#47) as row 48 on both the "All Trades" and
# "MarketMaking" logs - they mirror each other.
# ---------------------------------------------------------------------------
function Add-TradeRow48($ws) {
    $ws.Range("A48").Value = 47
    $ws.Range("D48").Value = "MarketMaking"
    $ws.Range("E48").Value = "DOWN"
    $ws.Range("F48").Value = 0.42
    $ws.Range("G48").Value = 0.41
    $ws.Range("H48").Value = "CLOSED"
    $ws.Range("I48").Value = -2.381
    $ws.Range("J48").Value = -0.01
    $ws.Range("K48").Value = 100.59
    $ws.Range("L48").Value = 0
    $ws.Range("M48").Value = 0
    $ws.Range("N48").Value = 0.6
    $ws.Range("O48").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P48").Value = "early_exit"
    $ws.Range("Q48").Value = 0.13

    # B48/C48 hold date-/time-look-alike text ("2026-02-17", "15:30:41").
    # A plain .Value assignment gets auto-coerced to a real date/time by
    # Excel's type inference, so instead build them as a text formula and
    # collapse the formula down to its literal text result via copy /
    # paste-special(values) - the same trick you'd use interactively to
    # keep a date-shaped string as text.
    $ws.Range("B48").Formula = '="2026-02-17"'
    $ws.Range("B48").Copy()
    $ws.Range("B48").PasteSpecial(-4163)

    $ws.Range("C48").Formula = '="15:30:41"'
    $ws.Range("C48").Copy()
    $ws.Range("C48").PasteSpecial(-4163)
}

Add-TradeRow48 $wb.Worksheets.Item("All Trades")
Add-TradeRow48 $wb.Worksheets.Item("MarketMaking")

$excel.CutCopyMode = $false

Write-Output "Applied trade #47 close-out updates."
